$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to
# Text format first, otherwise Excel auto-converts the string into a
# numeric value (losing formatting like trailing zeros).
$ws.Range('D2').Value = '29.809.77'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '1.897.18'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7631'
$ws.Range('E5').Value = '  +3.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '240.21'
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3040'
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '25.41'
$ws.Range('E9').Value = '  -3.51%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06824'
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07973'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7352'
$ws.Range('E12').Value = '  -4.53%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.859.10'
$ws.Range('E13').Value = '  -1.69%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.156'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '90.84'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '29.804.51'
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.80'
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.894'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '241.22'
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007696'
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '2.133.83'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.889'
$ws.Range('E24').Value = '  -0.90%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '166.55'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.211'
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.61'
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1285'
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.020'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.402'
$ws.Range('E30').Value = '  +3.43%  '
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.253'
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.060'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05193'
$ws.Range('E34').Value = '  +2.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.245'
$ws.Range('E35').Value = '  -2.30%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7237'
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.715'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.138'
$ws.Range('E40').Value = '  -2.53%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4389'
$ws.Range('E41').Value = '  -1.16%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '71.84'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.881'
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8273'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.596'
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '99.51'
$ws.Range('E47').Value = '  -1.45%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.753'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('D49').Value = '2.039.73'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '36.04'
$ws.Range('E50').Value = '  -0.88%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.471'
$ws.Range('E51').Value = '  +1.66%  '
